$ErrorActionPreference = "Stop"

$p = $ppt.ActivePresentation

$targetParaText = "Checkout de branch development"
$prefixText     = "Checkout de branch "
$boldText       = "development"

$found = $false

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)

        if (-not $shape.HasTextFrame) { continue }
        if (-not $shape.TextFrame.HasText) { continue }

        $tr = $shape.TextFrame.TextRange
        $paraCount = $tr.Paragraphs().Count

        for ($pi = 1; $pi -le $paraCount; $pi++) {
            $para = $tr.Paragraphs($pi, 1)
            $paraText = $para.Text.TrimEnd([char]13)

            if ($paraText -eq $targetParaText) {
                # Split the single run "Checkout de branch development" into
                # "Checkout de branch " (unchanged) + "development" (highlighted).
                $boldRange = $para.Characters($prefixText.Length + 1, $boldText.Length)

                $boldRange.Font.Bold = $true
                $boldRange.Font.Color.RGB = 65535          # RGB(255,255,0) -> FFFF00 yellow
                $boldRange.Font.Name = "Courier New"

                $found = $true
                break
            }
        }

        if ($found) { break }
    }

    if ($found) { break }
}

if (-not $found) {
    throw "Target paragraph 'Checkout de branch development' was not found."
}
